$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.320.21"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.14%  "

$ws.Range("D3").Value = "'1.869.13"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.60%  "

$ws.Range("D4").Value = "'0.9998"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.08%  "

$ws.Range("D5").Value = "'235.38"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.38%  "

$ws.Range("D6").Value = "'0.9999"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.08%  "

$ws.Range("D7").Value = "'0.4682"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.01%  "

$ws.Range("D8").Value = "'0.2845"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.72%  "

$ws.Range("D9").Value = "'0.06533"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.91%  "

$ws.Range("D10").Value = "'21.48"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.24%  "

$ws.Range("E11").Value = "  +1.40%  "

$ws.Range("D12").Value = "'97.94"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.60%  "

$ws.Range("D13").Value = "'1.865.18"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.77%  "

$ws.Range("E14").Value = "  +0.69%  "

$ws.Range("D15").Value = "'0.6772"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.67%  "

$ws.Range("D16").Value = "'279.47"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.24%  "

$ws.Range("D17").Value = "'30.308.15"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.07%  "

$ws.Range("D18").Value = "'0.9994"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.11%  "

$ws.Range("D19").Value = "'5.509"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.44%  "

$ws.Range("D20").Value = "'12.72"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.98%  "

$ws.Range("D21").Value = "'2.116.46"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.20%  "

$ws.Range("D22").Value = "'0.000007296"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.67%  "

$ws.Range("D23").Value = "'0.9997"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.09%  "

$ws.Range("D24").Value = "'6.167"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.07%  "

$ws.Range("D25").Value = "'9.192"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.61%  "

$ws.Range("D26").Value = "'165.23"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.53%  "

$ws.Range("D27").Value = "'19.14"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.05%  "

$ws.Range("D28").Value = "'1.932"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.55%  "

$ws.Range("D29").Value = "'1.378"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.39%  "

$ws.Range("D30").Value = "'0.09638"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.18%  "

$ws.Range("D31").Value = "'4.374"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.21%  "

$ws.Range("D32").Value = "'1.479"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.61%  "

$ws.Range("D33").Value = "'4.107"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.09%  "

$ws.Range("D34").Value = "'0.04714"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.07%  "

$ws.Range("D35").Value = "'1.130"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.38%  "

$ws.Range("D36").Value = "'0.7072"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.42%  "

$ws.Range("D37").Value = "'2.720"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.08%  "

$ws.Range("E38").Value = "  -0.46%  "

$ws.Range("D39").Value = "'6.276"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.02%  "

$ws.Range("D40").Value = "'2.533"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.59%  "

$ws.Range("D41").Value = "'73.82"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.68%  "

$ws.Range("D42").Value = "'1.949"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.36%  "

$ws.Range("D43").Value = "'0.8485"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.53%  "

$ws.Range("D44").Value = "'0.4183"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.20%  "

$ws.Range("D45").Value = "'103.96"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.96%  "

$ws.Range("E46").Value = "  -0.07%  "

$ws.Range("D47").Value = "'7.183"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.56%  "

$ws.Range("D48").Value = "'9.258"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.98%  "

$ws.Range("D49").Value = "'936.27"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.55%  "

$ws.Range("D50").Value = "'34.14"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.08%  "

$ws.Range("E51").Value = "  -1.76%  "
